# Updated cryptos list with latest Price (D) and Volume(1h) (E) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.008.99"
$ws.Range("E2").Value = "  -3.27%  "
$ws.Range("D3").Value = "1.714.17"
$ws.Range("E3").Value = "  -3.05%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'309.14"
$ws.Range("E5").Value = "  -5.92%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.4730"
$ws.Range("E7").Value = "  +4.24%  "
$ws.Range("D8").Value = "'0.3457"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").Value = "'41.90"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "'0.07231"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("D11").Value = "'1.037"
$ws.Range("E11").Value = "  -5.29%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "'19.76"
$ws.Range("E13").Value = "  -4.65%  "
$ws.Range("D14").Value = "'5.818"
$ws.Range("E14").Value = "  -3.21%  "
$ws.Range("D15").Value = "1.718.07"
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("D16").Value = "'6.822"
$ws.Range("E16").Value = "  -5.05%  "
$ws.Range("D17").Value = "'86.63"
$ws.Range("E17").Value = "  -6.44%  "
$ws.Range("D18").Value = "'0.00001035"
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("D19").Value = "'0.06378"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "'16.44"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").Value = "'5.603"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("D23").Value = "27.066.61"
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("E24").Value = "  -4.33%  "
$ws.Range("D25").Value = "'2.094"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "'19.89"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").Value = "'150.69"
$ws.Range("E27").Value = "  -5.38%  "
$ws.Range("D28").Value = "1.914.71"
$ws.Range("E28").Value = "  -3.04%  "
$ws.Range("D29").Value = "'2.053"
$ws.Range("E29").Value = "  -4.98%  "
$ws.Range("D30").Value = "'120.19"
$ws.Range("E30").Value = "  -3.24%  "
$ws.Range("D31").Value = "'1.022"
$ws.Range("E31").Value = "  -5.04%  "
$ws.Range("D32").Value = "'0.09132"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("D33").Value = "'3.599"
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("D34").Value = "'5.291"
$ws.Range("E34").Value = "  -5.79%  "
$ws.Range("D35").Value = "'1.468"
$ws.Range("E35").Value = "  +6.32%  "
$ws.Range("D36").Value = "'0.02172"
$ws.Range("E36").Value = "  -4.86%  "
$ws.Range("D37").Value = "'0.05832"
$ws.Range("E37").Value = "  -4.71%  "
$ws.Range("D38").Value = "'0.1996"
$ws.Range("E38").Value = "  -4.64%  "
$ws.Range("D39").Value = "'10.88"
$ws.Range("E39").Value = "  -8.12%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'4.692"
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("D42").Value = "'0.5952"
$ws.Range("E42").Value = "  -4.85%  "
$ws.Range("D43").Value = "'1.084"
$ws.Range("E43").Value = "  -8.04%  "
$ws.Range("D44").Value = "'7.460"
$ws.Range("E44").Value = "  -4.48%  "
$ws.Range("E45").Value = "  -3.62%  "
$ws.Range("D46").Value = "'3.575"
$ws.Range("E46").Value = "  -4.20%  "
$ws.Range("D47").Value = "'0.5557"
$ws.Range("E47").Value = "  -4.95%  "
$ws.Range("D48").Value = "'118.72"
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("D49").Value = "'1.819"
$ws.Range("E49").Value = "  -6.00%  "
$ws.Range("D50").Value = "'1.112"
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("D51").Value = "'0.06625"
